$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 39/40 full swap: PaxDollar/Quant reorder with updated price data
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D39") "105.68"
$ws.Range("E39").Value = "  +7.96%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D40") "0.9991"
$ws.Range("E40").Value = "  +0.02%  "

Set-TextValue $ws.Range("D2") "25.552.71"
$ws.Range("E2").Value = "  +2.63%  "
Set-TextValue $ws.Range("D3") "1.666.20"
$ws.Range("E3").Value = "  +1.79%  "
Set-TextValue $ws.Range("D4") "0.9983"
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "235.41"
$ws.Range("E5").Value = "  +1.56%  "
Set-TextValue $ws.Range("D6") "0.9996"
$ws.Range("E6").Value = "  -0.04%  "
Set-TextValue $ws.Range("D7") "0.4653"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("E8").Value = "  -0.13%  "
Set-TextValue $ws.Range("D9") "0.06137"
$ws.Range("E9").Value = "  +0.70%  "
Set-TextValue $ws.Range("D10") "1.663.68"
Set-TextValue $ws.Range("D11") "0.06962"
$ws.Range("E11").Value = "  -1.00%  "
Set-TextValue $ws.Range("D12") "14.67"
$ws.Range("E12").Value = "  +1.69%  "
Set-TextValue $ws.Range("D13") "4.351"
$ws.Range("E13").Value = "  +0.07%  "
Set-TextValue $ws.Range("D14") "74.90"
$ws.Range("E14").Value = "  +2.11%  "
Set-TextValue $ws.Range("D15") "0.5706"
$ws.Range("E15").Value = "  -3.98%  "
Set-TextValue $ws.Range("D16") "0.9998"
$ws.Range("E16").Value = "  +0.02%  "
Set-TextValue $ws.Range("D17") "0.9994"
$ws.Range("E17").Value = "  +0.03%  "
Set-TextValue $ws.Range("D18") "25.540.96"
$ws.Range("E18").Value = "  +2.58%  "
Set-TextValue $ws.Range("D19") "0.000006719"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("E20").Value = "  +1.53%  "
Set-TextValue $ws.Range("D21") "1.879.81"
$ws.Range("E21").Value = "  +1.74%  "
Set-TextValue $ws.Range("D22") "4.418"
$ws.Range("E22").Value = "  +1.42%  "
Set-TextValue $ws.Range("D23") "8.697"
$ws.Range("E23").Value = "  +1.84%  "
Set-TextValue $ws.Range("D24") "5.220"
$ws.Range("E24").Value = "  -0.21%  "
Set-TextValue $ws.Range("D25") "134.76"
$ws.Range("E25").Value = "  +1.11%  "
Set-TextValue $ws.Range("D26") "14.83"
$ws.Range("E26").Value = "  +0.17%  "
Set-TextValue $ws.Range("D27") "1.364"
$ws.Range("E27").Value = "  -1.31%  "
Set-TextValue $ws.Range("D28") "1.708"
$ws.Range("E28").Value = "  +4.66%  "
Set-TextValue $ws.Range("D29") "103.82"
$ws.Range("E29").Value = "  +0.35%  "
Set-TextValue $ws.Range("D30") "3.959"
$ws.Range("E30").Value = "  +3.13%  "
Set-TextValue $ws.Range("D31") "0.07711"
$ws.Range("E31").Value = "  +0.52%  "
Set-TextValue $ws.Range("D32") "3.604"
$ws.Range("E32").Value = "  +2.28%  "
Set-TextValue $ws.Range("D33") "0.04312"
$ws.Range("E33").Value = "  +0.68%  "
Set-TextValue $ws.Range("D34") "2.619"
$ws.Range("E34").Value = "  +1.71%  "
Set-TextValue $ws.Range("D35") "0.9442"
$ws.Range("E35").Value = "  +2.39%  "
Set-TextValue $ws.Range("D36") "0.5993"
$ws.Range("E36").Value = "  +2.87%  "
Set-TextValue $ws.Range("D37") "0.9176"
$ws.Range("E37").Value = "  +9.87%  "
Set-TextValue $ws.Range("D38") "2.477"
$ws.Range("E38").Value = "  -2.85%  "
Set-TextValue $ws.Range("D41") "0.01462"
$ws.Range("E41").Value = "  -3.51%  "
Set-TextValue $ws.Range("D42") "1.824"
$ws.Range("E42").Value = "  +5.64%  "
Set-TextValue $ws.Range("D43") "5.082"
$ws.Range("E43").Value = "  +8.78%  "
Set-TextValue $ws.Range("D44") "0.3710"
$ws.Range("E44").Value = "  +0.82%  "
Set-TextValue $ws.Range("D45") "0.1109"
$ws.Range("E45").Value = "  +2.86%  "
Set-TextValue $ws.Range("D46") "0.05251"
Set-TextValue $ws.Range("D47") "6.121"
$ws.Range("E47").Value = "  +1.52%  "
Set-TextValue $ws.Range("D48") "29.78"
$ws.Range("E48").Value = "  +2.17%  "
Set-TextValue $ws.Range("D49") "7.454"
$ws.Range("E49").Value = "  +3.58%  "
Set-TextValue $ws.Range("D50") "1.002"
$ws.Range("E50").Value = "  +0.21%  "
Set-TextValue $ws.Range("D51") "0.9986"
$ws.Range("E51").Value = "  +0.17%  "
